# Replace the quarter-label text in column A (rows 2-18) with real dates
# (quarter start dates), formatted as short dates, matching the "added
# updated data and plots" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> quarter start date (these replace shared-string labels
# like "2020Q2", "Q3", "Q4", "Q1" with actual date values).
$dates = @{
    2  = "4/1/2020"
    3  = "7/1/2020"
    4  = "10/1/2020"
    5  = "1/1/2020"
    6  = "4/1/2021"
    7  = "7/1/2021"
    8  = "10/1/2021"
    9  = "1/1/2021"
    10 = "4/1/2022"
    11 = "7/1/2022"
    12 = "10/1/2022"
    13 = "1/1/2022"
    14 = "4/1/2023"
    15 = "7/1/2023"
    16 = "10/1/2023"
    17 = "1/1/2023"
    18 = "1/1/2024"
}

foreach ($row in 2..18) {
    $cell = $ws.Cells.Item($row, 1)
    # Set the number format before the value so the cell doesn't pick up
    # a transient "auto-detected" date format in addition to this one.
    $cell.NumberFormat = "mm-dd-yy"
    $cell.Value = $dates[[string]$row]
}

# Column A is now date-valued; widen/best-fit it like Excel would after
# the data was pasted in.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Update the view: zoom in on the sheet and move the active selection.
$window = $excel.ActiveWindow
$window.Zoom = 128
$ws.Range("A19").Select() | Out-Null
